$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.059.82"
$ws.Range("E2").Value = "  +0.56%  "
$ws.Range("D3").Value = "1.678.71"
$ws.Range("E3").Value = "  +0.64%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.65"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.08%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.517"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -1.17%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.255"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +2.39%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "21.41"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +5.26%  "
$ws.Range("E10").Value = "  +0.27%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0889"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -0.60%  "
$ws.Range("D12").Value = "1.915.14"
$ws.Range("E12").Value = "  +0.58%  "
$ws.Range("D13").Value = "1.664.72"
$ws.Range("E13").Value = "  -0.11%  "
$ws.Range("E14").Value = "  +0.92%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.535"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +1.51%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "66.42"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +0.47%  "
$ws.Range("D17").Value = "27.058.02"
$ws.Range("E17").Value = "  +0.56%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "8.17"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +2.21%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "235.87"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +0.56%  "
$ws.Range("D20").Value = "0.0₃0736"
$ws.Range("E20").Value = "  +0.57%  "
$ws.Range("E21").Value = "  -0.02%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.47"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +1.57%  "
$ws.Range("E23").Value = "  +1.38%  "
$ws.Range("E24").Value = "  -2.68%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "147.32"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.46%  "
$ws.Range("E26").Value = "  +1.99%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.53"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +4.00%  "
$ws.Range("E29").Value = "  +0.21%  "
$ws.Range("E30").Value = "  +0.68%  "
$ws.Range("E31").Value = "  +0.03%  "
$ws.Range("E32").Value = "  +1.05%  "
$ws.Range("D33").Value = "1.539.16"
$ws.Range("E33").Value = "  +6.13%  "
$ws.Range("E34").Value = "  +1.15%  "
$ws.Range("E35").Value = "  +4.83%  "
$ws.Range("E36").Value = "  -1.01%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.588"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +0.62%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.916"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +1.48%  "
$ws.Range("E39").Value = "  +2.27%  "
$ws.Range("E40").Value = "  +6.07%  "
$ws.Range("E41").Value = "  +0.10%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "68.02"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +3.04%  "
$ws.Range("E43").Value = "  -3.67%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.27"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -1.26%  "
$ws.Range("D45").Value = "1.822.94"
$ws.Range("E45").Value = "  +0.43%  "
$ws.Range("E46").Value = "  -0.19%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "90.37"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -0.51%  "
$ws.Range("B48").Value = "RenderToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.54"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -0.02%  "
$ws.Range("B49").Value = "Algorand"
$ws.Range("C49").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.104"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +1.76%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.04"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +6.63%  "
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0507"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +0.10%  "
